# Insert two new weekly price rows for Pimiento ("Zafiro rojo" / "Zafiro verde")
# at the top of the existing 337-374 block, pushing the rest of the block
# (and the sheet's used range) down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 337; this shifts old rows
# 337..374 down to 339..376 and carries the date-format style (s="2")
# that column D already had at that position.
$ws.Rows.Item(337).EntireRow.Insert()
$ws.Rows.Item(337).EntireRow.Insert()

# --- New row 337 ---
$ws.Range("A337").Value = 7
$ws.Range("B337").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C337").Value = "Ñuble"
$ws.Range("D337").Value = 44918
$ws.Range("E337").Value = 16
$ws.Range("F337").Value = 100112002
$ws.Range("G337").Value = "Pimiento"
$ws.Range("H337").Value = "Zafiro rojo"
$ws.Range("I337").Value = "Primera"
$ws.Range("J337").Value = 100
$ws.Range("K337").Value = 12000
$ws.Range("L337").Value = 13000
$ws.Range("M337").Value = 12500
$ws.Range("N337").Value = "$/caja 15 kilos"
$ws.Range("O337").Value = "Región de Arica y Parinacota"
$ws.Range("P337").Value = 833
$ws.Range("Q337").Value = 15
$ws.Range("R337").Value = "Hortaliza"

# --- New row 338 ---
$ws.Range("A338").Value = 7
$ws.Range("B338").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C338").Value = "Ñuble"
$ws.Range("D338").Value = 44918
$ws.Range("E338").Value = 16
$ws.Range("F338").Value = 100112002
$ws.Range("G338").Value = "Pimiento"
$ws.Range("H338").Value = "Zafiro verde"
$ws.Range("I338").Value = "Primera"
$ws.Range("J338").Value = 100
$ws.Range("K338").Value = 11000
$ws.Range("L338").Value = 12000
$ws.Range("M338").Value = 11500
$ws.Range("N338").Value = "$/caja 15 kilos"
$ws.Range("O338").Value = "Región de Arica y Parinacota"
$ws.Range("P338").Value = 767
$ws.Range("Q338").Value = 15
$ws.Range("R338").Value = "Hortaliza"
